$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new booking row (row 3) below the existing data
$ws.Range("A3").Value = "Akhila"
$ws.Range("B3").Value = "Batchu"
$ws.Range("C3").Value = 345
$ws.Range("D3").Value = $false
$ws.Range("E3").Value = 44682
$ws.Range("F3").Value = 45078
$ws.Range("G3").Value = "lunch"

# Match the date formatting used by the row above (reuse its number format)
$ws.Range("E2:F2").Copy()
$ws.Range("E3:F3").PasteSpecial(-4122)

# Leave the selection where the user ended up after entering the data
$ws.Range("C5").Select()
